# Update the "want to go" counters (column F) for several rows across
# multiple worksheets, per the scraped-data refresh described in the diff.

$wb = $excel.ActiveWorkbook

# Map: worksheet name -> list of (row, newValue)
$updates = @{
    "展览"     = @{ 2 = 680; 6 = 604; 7 = 29; 8 = 346; 10 = 6079; 12 = 1065; 13 = 26; 17 = 595; 18 = 1022; 19 = 53; 20 = 37; 21 = 195; 22 = 1353; 25 = 74 }
    "演出"     = @{ 15 = 84; 16 = 632 }
    "本地生活" = @{ 10 = 121 }
    "全部类型" = @{ 9 = 121; 10 = 680; 13 = 604; 15 = 29; 16 = 346; 18 = 6079; 21 = 1065; 26 = 595; 28 = 84; 31 = 1022; 32 = 37; 40 = 74 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
